$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "post a money transfer" requirement (row 10) as DONE
$ws.Range("C10").Value = "DONE"

# Update the active selection to C11 as recorded in the saved workbook
$ws.Range("C11").Select()
